# Add a new Job Posting row with Job_Id = JD_011
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$jobDescription = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

# New row goes right after the last existing data row (row 11 -> row 12)
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "JD_011"
$ws.Cells.Item($newRow, 2).Value = "Cyber Security Engineer"
$ws.Cells.Item($newRow, 3).Value = $jobDescription
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 4
